# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps produced by the handback report
# generator, without touching anything else in the workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first row.
$wsOverview.Range("G2").Value = "2016-09-09 12:52:43"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime".
$wsZhCn.Range("H2").Value = "2016-09-09 12:52:32"
$wsZhCn.Range("K2").Value = "2016-09-09 12:53:29"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime".
$wsDeDe.Range("H2").Value = "2016-09-09 12:52:43"
$wsDeDe.Range("K2").Value = "2016-09-09 12:53:47"
